# Refined metadata to be additional tab
#
# 1) Add a new "metadata" worksheet after the existing "data" sheet.
# 2) Populate it with the panel-level metadata (headers + one data row).
# 3) Refresh the per-row "time_taken" timestamps on the "data" sheet to the
#    new query time.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. New "metadata" sheet, placed right after "data"
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Add($null, $dataSheet)
$meta.Name = "metadata"

# ---------------------------------------------------------------------
# 2. Header row (bold, bordered, centered/top-aligned - same look as the
#    header row already used on the "data" sheet). Re-use the existing
#    header style by copying it over, instead of re-building it property
#    by property (keeps the same style entry instead of cloning a new one).
# ---------------------------------------------------------------------
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$dataSheet.Range("B1:F1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)

$meta.Range("A2").Value = 0
$dataSheet.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Data row - keep the version string as text ("1.56"), not a number
# ---------------------------------------------------------------------
$meta.Range("B2").Value = "Rhabdomyolysis and metabolic muscle disorders"
$meta.Range("C2").Value = 66

$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.56"
$meta.Range("D2").ClearFormats()

$meta.Range("E2").Value = "2021-09-02T16:36:28.455913Z"
$meta.Range("F2").Value = "2021-10-05 14:22:36.375695"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/66/?format=json"

[void]$meta.Range("A1").Select()

# ---------------------------------------------------------------------
# 3. Refresh "time_taken" column on the "data" sheet
# ---------------------------------------------------------------------
$newTimes = @(
    "2021-10-05 14:22:36.379165",
    "2021-10-05 14:22:36.379173",
    "2021-10-05 14:22:36.379176",
    "2021-10-05 14:22:36.379179",
    "2021-10-05 14:22:36.379181",
    "2021-10-05 14:22:36.379184",
    "2021-10-05 14:22:36.379186",
    "2021-10-05 14:22:36.379189",
    "2021-10-05 14:22:36.379192",
    "2021-10-05 14:22:36.379194",
    "2021-10-05 14:22:36.379196",
    "2021-10-05 14:22:36.379199",
    "2021-10-05 14:22:36.379201",
    "2021-10-05 14:22:36.379204",
    "2021-10-05 14:22:36.379206",
    "2021-10-05 14:22:36.379209",
    "2021-10-05 14:22:36.379211",
    "2021-10-05 14:22:36.379214",
    "2021-10-05 14:22:36.379216",
    "2021-10-05 14:22:36.379219",
    "2021-10-05 14:22:36.379221",
    "2021-10-05 14:22:36.379223",
    "2021-10-05 14:22:36.379226",
    "2021-10-05 14:22:36.379229",
    "2021-10-05 14:22:36.379231",
    "2021-10-05 14:22:36.379234",
    "2021-10-05 14:22:36.379236",
    "2021-10-05 14:22:36.379239",
    "2021-10-05 14:22:36.379241",
    "2021-10-05 14:22:36.379244",
    "2021-10-05 14:22:36.379246",
    "2021-10-05 14:22:36.379249",
    "2021-10-05 14:22:36.379251",
    "2021-10-05 14:22:36.379254",
    "2021-10-05 14:22:36.379256",
    "2021-10-05 14:22:36.379259",
    "2021-10-05 14:22:36.379261",
    "2021-10-05 14:22:36.379264",
    "2021-10-05 14:22:36.379266",
    "2021-10-05 14:22:36.379269",
    "2021-10-05 14:22:36.379271",
    "2021-10-05 14:22:36.379274",
    "2021-10-05 14:22:36.379276",
    "2021-10-05 14:22:36.379279",
    "2021-10-05 14:22:36.379281",
    "2021-10-05 14:22:36.379284",
    "2021-10-05 14:22:36.379286",
    "2021-10-05 14:22:36.379288",
    "2021-10-05 14:22:36.379291",
    "2021-10-05 14:22:36.379293",
    "2021-10-05 14:22:36.379296",
    "2021-10-05 14:22:36.379298",
    "2021-10-05 14:22:36.379301",
    "2021-10-05 14:22:36.379304",
    "2021-10-05 14:22:36.379306",
    "2021-10-05 14:22:36.379308",
    "2021-10-05 14:22:36.379311",
    "2021-10-05 14:22:36.379313",
    "2021-10-05 14:22:36.379316",
    "2021-10-05 14:22:36.379318"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $newTimes[$i]
}
